$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row of data ---
$logs = $wb.Worksheets.Item("Logs")

$newRow = 31

$logs.Cells.Item($newRow, 1).Value = "Wanneer zijn jullie open?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Cells.Item($newRow, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,
Bedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.
Met vriendelijke groet,
[Naam bedrijf]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-26 23:03:24"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"

# --- Extend the conditional formatting ranges to cover the new row ---
$ranges = @("D2:D30", "G2:G30", "H2:H30", "I2:I30")
$newRanges = @("D2:D31", "G2:G31", "H2:H31", "I2:I31")

for ($i = 0; $i -lt $ranges.Length; $i++) {
    $fcs = $logs.Range($ranges[$i]).FormatConditions
    $target = $logs.Range($newRanges[$i])
    for ($j = 1; $j -le $fcs.Count; $j++) {
        $fcs.Item($j).ModifyAppliesToRange($target)
    }
}

# --- Sheet "Dashboard": increment the count for "Openingstijden / Locatie" ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(3, 2).Value = 8
